# Update crypto price/volume figures per the GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.37"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.41%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.106"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.42%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08209"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.44%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.073"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.48%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.950"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.23%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "9.29%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9246"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.09%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1120"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.26%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1909"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.45%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09158"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.74%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03650"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.33%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09918"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.14%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001435"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.38%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005851"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.82%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.481"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.00%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.127"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.33%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3425"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.05%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1309"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.23%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.089"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.59%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.94%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.73%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.58%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004810"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.29%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001249"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.84%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004447"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.16%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01971"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.92%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04871"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.06%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007667"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.13%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009175"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "18.54%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.97%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002084"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.19%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01162"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.36%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006549"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.90%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "179.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "247.94%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-21.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
